# PowerShell/COM script to refresh Sagittarius_Profits leve-profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) with the
# latest market-board pricing pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Worksheet "ALC" ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21: Book and a Hard Place
$ws.Cells.Item(21, 8).Value = 14999
$ws.Cells.Item(21, 9).Value = 14999
$ws.Cells.Item(21, 11).Value = 14999
$ws.Cells.Item(21, 13).Value = -14531
# Row 23: There's Something about Bury
$ws.Cells.Item(23, 8).Value = 14999
$ws.Cells.Item(23, 9).Value = 14999
$ws.Cells.Item(23, 11).Value = 14999
$ws.Cells.Item(23, 13).Value = -14765
# Row 32: Automata for the People
$ws.Cells.Item(32, 8).Value = 3434
$ws.Cells.Item(32, 9).Value = 3000
$ws.Cells.Item(32, 10).Value = 3607.6
$ws.Cells.Item(32, 11).Value = 3000
$ws.Cells.Item(32, 12).Value = 3607.6
$ws.Cells.Item(32, 13).Value = -2674
$ws.Cells.Item(32, 14).Value = -4259.6
# Row 62: The Mustache Suits Him
$ws.Cells.Item(62, 8).Value = 5049.1665
$ws.Cells.Item(62, 9).Value = 3573.75
$ws.Cells.Item(62, 11).Value = 3573.75
$ws.Cells.Item(62, 13).Value = -2949.75
# Row 65: Forgery of Convenience (L)
$ws.Cells.Item(65, 8).Value = 5049.1665
$ws.Cells.Item(65, 9).Value = 3573.75
$ws.Cells.Item(65, 11).Value = 17868.75
$ws.Cells.Item(65, 13).Value = -14748.75
# Row 70: Consecrating Congregation
$ws.Cells.Item(70, 8).Value = 70472.87
$ws.Cells.Item(70, 10).Value = 75342.36
$ws.Cells.Item(70, 12).Value = 226027.08
$ws.Cells.Item(70, 14).Value = -226567.08
# Row 73: Curbing the Contagion (L)
$ws.Cells.Item(73, 8).Value = 70472.87
$ws.Cells.Item(73, 10).Value = 75342.36
$ws.Cells.Item(73, 12).Value = 226027.08
$ws.Cells.Item(73, 14).Value = -227899.08
# Row 125: Body over Mind
$ws.Cells.Item(125, 8).Value = 250005740
$ws.Cells.Item(125, 9).Value = 250008480
$ws.Cells.Item(125, 10).Value = 250000270
$ws.Cells.Item(125, 11).Value = 2250076320
$ws.Cells.Item(125, 12).Value = 2250002430
$ws.Cells.Item(125, 13).Value = -2250073860
$ws.Cells.Item(125, 14).Value = -2250007350
# Row 131: Mindful Study
$ws.Cells.Item(131, 8).Value = 12548.75
$ws.Cells.Item(131, 9).Value = 3947
$ws.Cells.Item(131, 10).Value = 15416
$ws.Cells.Item(131, 11).Value = 11841
$ws.Cells.Item(131, 12).Value = 46248
$ws.Cells.Item(131, 13).Value = -6801
$ws.Cells.Item(131, 14).Value = -56328
# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 1357.8334
$ws.Cells.Item(137, 10).Value = 1383
$ws.Cells.Item(137, 12).Value = 4149
$ws.Cells.Item(137, 14).Value = -9249
# Row 141: Remedy for Reason
$ws.Cells.Item(141, 8).Value = 3699.4
$ws.Cells.Item(141, 9).Value = 4166.6665
$ws.Cells.Item(141, 10).Value = 2998.5
$ws.Cells.Item(141, 11).Value = 12499.9995
$ws.Cells.Item(141, 12).Value = 8995.5
$ws.Cells.Item(141, 13).Value = -7319.999500000002
$ws.Cells.Item(141, 14).Value = -19355.5

# --- Worksheet "ARM" ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 4002.7856
$ws.Cells.Item(61, 9).Value = 866.4286
$ws.Cells.Item(61, 11).Value = 866.4286
$ws.Cells.Item(61, 13).Value = -654.4286
# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 1144.4546
$ws.Cells.Item(132, 9).Value = 1144.4546
$ws.Cells.Item(132, 11).Value = 3433.3638
$ws.Cells.Item(132, 13).Value = -903.3638000000001
# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 4002.7856
$ws.Cells.Item(136, 9).Value = 866.4286
$ws.Cells.Item(136, 11).Value = 2599.2858
$ws.Cells.Item(136, 13).Value = -49.28579999999965

# --- Worksheet "BSM" ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Cells.Item(22, 8).Value = 890.7143
$ws.Cells.Item(22, 9).Value = 649
$ws.Cells.Item(22, 10).Value = 1213
$ws.Cells.Item(22, 11).Value = 649
$ws.Cells.Item(22, 12).Value = 1213
$ws.Cells.Item(22, 13).Value = -476
$ws.Cells.Item(22, 14).Value = -1559
# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 1597.25
$ws.Cells.Item(86, 9).Value = 1706
$ws.Cells.Item(86, 11).Value = 1706
$ws.Cells.Item(86, 13).Value = -583
# Row 88: Swords for Plowshares
$ws.Cells.Item(88, 8).Value = 33397.25
$ws.Cells.Item(88, 10).Value = 33397.25
$ws.Cells.Item(88, 12).Value = 33397.25
$ws.Cells.Item(88, 14).Value = -34209.25
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 1597.25
$ws.Cells.Item(89, 9).Value = 1706
$ws.Cells.Item(89, 11).Value = 8530
$ws.Cells.Item(89, 13).Value = -2914
# Row 91: Negative, They Are Meat Popsicles (L)
$ws.Cells.Item(91, 8).Value = 33397.25
$ws.Cells.Item(91, 10).Value = 33397.25
$ws.Cells.Item(91, 12).Value = 33397.25
$ws.Cells.Item(91, 14).Value = -36205.25
# Row 105: Ingot to Wing It
$ws.Cells.Item(105, 8).Value = 2090.6316
$ws.Cells.Item(105, 9).Value = 1920.1875
$ws.Cells.Item(105, 11).Value = 1920.1875
$ws.Cells.Item(105, 13).Value = -173.1875
# Row 107: The Gold Experience
$ws.Cells.Item(107, 8).Value = 2261.8
$ws.Cells.Item(107, 9).Value = 2261.8
$ws.Cells.Item(107, 11).Value = 2261.8
$ws.Cells.Item(107, 13).Value = -341.8000000000002

# --- Worksheet "CRP" ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 2646.75
$ws.Cells.Item(58, 10).Value = 2533
$ws.Cells.Item(58, 12).Value = 2533
$ws.Cells.Item(58, 14).Value = -2939
# Row 59: Bow Down to Magic
$ws.Cells.Item(59, 8).Value = 117333
$ws.Cells.Item(59, 9).Value = 26000
$ws.Cells.Item(59, 11).Value = 26000
$ws.Cells.Item(59, 13).Value = -24855
# Row 60: Bowing to Greater Power
$ws.Cells.Item(60, 8).Value = 14599.667
$ws.Cells.Item(60, 9).Value = 4650
$ws.Cells.Item(60, 11).Value = 4650
$ws.Cells.Item(60, 13).Value = -4139
# Row 105: Zelkova, My Love
$ws.Cells.Item(105, 8).Value = 3138.077
$ws.Cells.Item(105, 9).Value = 2554.818
$ws.Cells.Item(105, 10).Value = 3565.8
$ws.Cells.Item(105, 11).Value = 2554.818
$ws.Cells.Item(105, 12).Value = 3565.8
$ws.Cells.Item(105, 13).Value = -807.8180000000002
$ws.Cells.Item(105, 14).Value = -7059.8
# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 2084.375
$ws.Cells.Item(132, 10).Value = 2266.3333
$ws.Cells.Item(132, 12).Value = 6798.999899999999
$ws.Cells.Item(132, 14).Value = -11858.9999
# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 2646.75
$ws.Cells.Item(136, 10).Value = 2533
$ws.Cells.Item(136, 12).Value = 7599
$ws.Cells.Item(136, 14).Value = -12699

# --- Worksheet "CUL" ---
$ws = $wb.Worksheets.Item("CUL")
# Row 56: Culture Club
$ws.Cells.Item(56, 8).Value = 12060
$ws.Cells.Item(56, 9).Value = 12060
$ws.Cells.Item(56, 11).Value = 12060
$ws.Cells.Item(56, 13).Value = -11530
# Row 119: Super Dark Times
$ws.Cells.Item(119, 8).Value = 2963.3333
$ws.Cells.Item(119, 9).Value = 2963.3333
$ws.Cells.Item(119, 11).Value = 8889.999899999999
$ws.Cells.Item(119, 13).Value = -4051.999899999999
# Row 138: Bring Me Your Tacos
$ws.Cells.Item(138, 8).Value = 1430887.9
$ws.Cells.Item(138, 9).Value = 1430887.9
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 4292663.699999999
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 13).Value = -4287523.699999999
$ws.Cells.Item(138, 14).ClearContents()
# Row 141: Ocean Explosion
$ws.Cells.Item(141, 8).Value = 7642.5713
$ws.Cells.Item(141, 9).Value = 7642.5713
$ws.Cells.Item(141, 11).Value = 22927.7139
$ws.Cells.Item(141, 13).Value = -17747.7139

# --- Worksheet "GSM" ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Cells.Item(70, 8).Value = 6894.6665
$ws.Cells.Item(70, 9).Value = 6795.4116
$ws.Cells.Item(70, 10).Value = 7000.125
$ws.Cells.Item(70, 11).Value = 6795.4116
$ws.Cells.Item(70, 12).Value = 7000.125
$ws.Cells.Item(70, 13).Value = -6525.4116
$ws.Cells.Item(70, 14).Value = -7540.125
# Row 73: Hulls of Broken Dreams (L)
$ws.Cells.Item(73, 8).Value = 6894.6665
$ws.Cells.Item(73, 9).Value = 6795.4116
$ws.Cells.Item(73, 10).Value = 7000.125
$ws.Cells.Item(73, 11).Value = 6795.4116
$ws.Cells.Item(73, 12).Value = 7000.125
$ws.Cells.Item(73, 13).Value = -5859.4116
$ws.Cells.Item(73, 14).Value = -8872.125
# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 1559.7273
$ws.Cells.Item(102, 9).Value = 1465.7
$ws.Cells.Item(102, 11).Value = 1465.7
$ws.Cells.Item(102, 13).Value = 156.3
# Row 113: Copious Crystal Cannons
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 1569
$ws.Cells.Item(122, 9).Value = 1556.4546
$ws.Cells.Item(122, 10).Value = 1603.5
$ws.Cells.Item(122, 11).Value = 4669.3638
$ws.Cells.Item(122, 12).Value = 4810.5
$ws.Cells.Item(122, 13).Value = -2219.3638
$ws.Cells.Item(122, 14).Value = -9710.5
# Row 126: Gold Rush Order
$ws.Cells.Item(126, 8).Value = 3485.375
$ws.Cells.Item(126, 9).Value = 3549
$ws.Cells.Item(126, 11).Value = 10647
$ws.Cells.Item(126, 13).Value = -8177

# --- Worksheet "LTW" ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Cells.Item(40, 8).Value = 2651.12
$ws.Cells.Item(40, 9).Value = 2129.4707
$ws.Cells.Item(40, 11).Value = 2129.4707
$ws.Cells.Item(40, 13).Value = -1993.4707
# Row 82: Trainin' the Neck
$ws.Cells.Item(82, 8).Value = 1821.1052
$ws.Cells.Item(82, 9).Value = 2487.5833
$ws.Cells.Item(82, 10).Value = 678.5714
$ws.Cells.Item(82, 11).Value = 2487.5833
$ws.Cells.Item(82, 12).Value = 678.5714
$ws.Cells.Item(82, 13).Value = -2126.5833
$ws.Cells.Item(82, 14).Value = -1400.5714
# Row 85: Training Is Only Skintight (L)
$ws.Cells.Item(85, 8).Value = 1821.1052
$ws.Cells.Item(85, 9).Value = 2487.5833
$ws.Cells.Item(85, 10).Value = 678.5714
$ws.Cells.Item(85, 11).Value = 2487.5833
$ws.Cells.Item(85, 12).Value = 678.5714
$ws.Cells.Item(85, 13).Value = -1239.5833
$ws.Cells.Item(85, 14).Value = -3174.5714
# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 3999.5
$ws.Cells.Item(136, 9).Value = 3000
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 13).Value = -6450

# --- Worksheet "WVR" ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Cells.Item(81, 8).Value = 1252652.4
$ws.Cells.Item(81, 9).Value = 4199.5
$ws.Cells.Item(81, 10).Value = 2501105.2
$ws.Cells.Item(81, 11).Value = 8399
$ws.Cells.Item(81, 12).Value = 5002210.4
$ws.Cells.Item(81, 13).Value = -7338
$ws.Cells.Item(81, 14).Value = -5004332.4
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Cells.Item(84, 8).Value = 1252652.4
$ws.Cells.Item(84, 9).Value = 4199.5
$ws.Cells.Item(84, 10).Value = 2501105.2
$ws.Cells.Item(84, 11).Value = 41995
$ws.Cells.Item(84, 12).Value = 25011052
$ws.Cells.Item(84, 13).Value = -36691
$ws.Cells.Item(84, 14).Value = -25021660
# Row 100: Of Great Import
$ws.Cells.Item(100, 8).Value = 5556675
$ws.Cells.Item(100, 9).Value = 6251072
$ws.Cells.Item(100, 11).Value = 12502144
$ws.Cells.Item(100, 13).Value = -12501603
# Row 113: A Tender Table
$ws.Cells.Item(113, 8).Value = 1999.5
$ws.Cells.Item(113, 10).Value = 1999
$ws.Cells.Item(113, 12).Value = 5997
$ws.Cells.Item(113, 14).Value = -10337
# Row 122: Heavy Armoire
$ws.Cells.Item(122, 8).Value = 1280.6
$ws.Cells.Item(122, 9).Value = 1280.6
$ws.Cells.Item(122, 11).Value = 3841.8
$ws.Cells.Item(122, 13).Value = -1391.8
# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 3061.2727
$ws.Cells.Item(132, 9).Value = 4211.7856
$ws.Cells.Item(132, 10).Value = 1047.875
$ws.Cells.Item(132, 11).Value = 12635.3568
$ws.Cells.Item(132, 12).Value = 3143.625
$ws.Cells.Item(132, 13).Value = -10105.3568
$ws.Cells.Item(132, 14).Value = -8203.625
# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 2790.7742
$ws.Cells.Item(136, 9).Value = 2821.5417
$ws.Cells.Item(136, 10).Value = 2685.2856
$ws.Cells.Item(136, 11).Value = 8464.625100000001
$ws.Cells.Item(136, 12).Value = 8055.8568
$ws.Cells.Item(136, 13).Value = -5914.625100000001
$ws.Cells.Item(136, 14).Value = -13155.8568
